# Commit: "remove 4 new vars for now"
# Removes the 4 newly-added columns (Site ID, Tank #, Calibration Date,
# Recheck Year) -- both their header/sub-header labels (row 17/18 and the
# equivalent row 43/44 block below) and the associated Carbone template
# placeholders, which were stored in columns H:K.
# Clearing the cell contents also drops the now-unused shared strings
# automatically when the workbook is re-saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H17:K17").ClearContents() | Out-Null
$ws.Range("H18:K18").ClearContents() | Out-Null
$ws.Range("H43:K43").ClearContents() | Out-Null
$ws.Range("H44:K44").ClearContents() | Out-Null

# Restore the selection/cursor position that results from the edit
# (previously scrolled to A15 with F50 selected; after the edit the
# view is back at the top with I38 selected).
$ws.Range("I38").Select() | Out-Null
